# Semana 33 de 2025 - update poisson.xlsx statistics table (rows 2-35)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a string value while preserving it as Text
# (needed for numeric-looking codes stored in column A, e.g. "330")
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# --- Column A event codes (numeric-looking, must stay text) ---
Set-TextValue "A10" "330"
Set-TextValue "A11" "340"
Set-TextValue "A12" "342"
Set-TextValue "A13" "346"
Set-TextValue "A14" "348"
Set-TextValue "A15" "352"
Set-TextValue "A16" "355"
Set-TextValue "A17" "356"
Set-TextValue "A18" "357"
Set-TextValue "A19" "365"
Set-TextValue "A20" "420"
Set-TextValue "A21" "455"
Set-TextValue "A22" "465"
Set-TextValue "A23" "535"
Set-TextValue "A24" "549"
Set-TextValue "A25" "560"
Set-TextValue "A26" "580"
Set-TextValue "A27" "591"
Set-TextValue "A28" "610"

# --- Column B event names ---
$ws.Range("B10").Value = "Hepatitis a"
$ws.Range("B11").Value = "Hepatitis b, c y coinfeccion hepatitis b y delta"
$ws.Range("B12").Value = "Enfermedades huerfanas - raras"
$ws.Range("B13").Value = "Ira por virus nuevo"
$ws.Range("B14").Value = "Infeccion respiratoria aguda grave irag inusitada"
$ws.Range("B15").Value = "Infecciones de sitio quirurgico asociadas a procedimiento medico quirurgico"
$ws.Range("B16").Value = "Enfermedad transmitida por alimentos o agua (eta)"
$ws.Range("B17").Value = "Intento de suicidio"
$ws.Range("B18").Value = "Iad - infecciones asociadas a dispositivos - individual"
$ws.Range("B19").Value = "Intoxicaciones"
$ws.Range("B20").Value = "Leishmaniasis cutanea"
$ws.Range("B21").Value = "Leptospirosis"
$ws.Range("B22").Value = "Malaria"
$ws.Range("B23").Value = "Meningitis bacteriana y enfermedad meningoc”cica"
$ws.Range("B24").Value = "Morbilidad materna extrema"
$ws.Range("B25").Value = "Mortalidad perinatal y neonatal tardia"
$ws.Range("B26").Value = "Mortalidad por dengue"
$ws.Range("B27").Value = "Vigilancia integrada de muertes en menores de cinco anos por infeccion respiratoria aguda - enfermedad diarreica aguda y/o desnutricion"

# --- Numeric columns C (Esperado), D (Observado), E (valor p) ---
$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 0
$ws.Range("D4").Value = 1
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 11
$ws.Range("E5").Value = 0.02
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 0.05
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 0
$ws.Range("C9").Value = 46
$ws.Range("D9").Value = 34
$ws.Range("E9").Value = 0.01
$ws.Range("C10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 0
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0.15
$ws.Range("C13").Value = 2
$ws.Range("E13").Value = 0.14
$ws.Range("C14").Value = 1
$ws.Range("E14").Value = 0.37
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 1
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 0.14
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 0.12
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0.14
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0.12
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 1
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 0.27
$ws.Range("C22").Value = 0
$ws.Range("E22").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = 0.15
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0.37
$ws.Range("C26").Value = 0
$ws.Range("E26").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0.37
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = 0.18
$ws.Range("C33").Value = 10
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = 0
$ws.Range("C34").Value = 8
$ws.Range("D34").Value = 3
$ws.Range("E34").Value = 0.03
$ws.Range("C35").Value = 10
$ws.Range("D35").Value = 7
$ws.Range("E35").Value = 0.09

# --- Clear cells that no longer have data (row 28 now mirrors former row 10 shape) ---
$ws.Range("B28").Value = ""
$ws.Range("C28").Value = ""
$ws.Range("E28").Value = ""

